# #5: insurance, claim, debt, investment done
#
# The "債務" (debt) sheet gains new columns: species, debtor, total are
# inserted into the header row (between the existing owner/register_date
# columns), and the row is extended out to match the same
# property_category..index tail that the "汽車" sheet already has. Row 2
# picks up matching data, including a brand new "debt" value in the new
# H column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# ---- Header row (row 1) -------------------------------------------------
# Existing: B1=species-slot(was "中期放款"/21), C1=owner(was 14), D1=register_reason(was 22)...
# New layout walks left-to-right so brand-new shared strings ("species",
# "debtor", "total") get interned in the same order the diff expects.
$ws.Cells.Item(1, 2).Value = "species"            # B1 (new)
$ws.Cells.Item(1, 3).Value = "debtor"              # C1 (new)
$ws.Cells.Item(1, 4).Value = "owner"               # D1 (existing string)
$ws.Cells.Item(1, 5).Value = "total"               # E1 (new)
$ws.Cells.Item(1, 6).Value = "register_date"       # F1 (existing string)
$ws.Cells.Item(1, 7).Value = "register_reason"     # G1 (existing string)
$ws.Cells.Item(1, 8).Value = "property_category"   # H1 (new cell)
$ws.Cells.Item(1, 9).Value = "category"            # I1 (new cell)
$ws.Cells.Item(1, 10).Value = "date"                # J1 (new cell)
$ws.Cells.Item(1, 11).Value = "legislator_name"     # K1 (new cell)
$ws.Cells.Item(1, 12).Value = "legislator_id"       # L1 (new cell)
$ws.Cells.Item(1, 13).Value = "source_file"         # M1 (new cell)
$ws.Cells.Item(1, 14).Value = "index"               # N1 (new cell)

# Give the newly added header cells (H1:N1) the same bold/border style as
# the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# ---- Data row (row 2) ---------------------------------------------------
$ws.Cells.Item(2, 2).Value = "中期放款"              # B2 (species)
$ws.Cells.Item(2, 3).Value = "陳亭妃"                # C2 (debtor)
$ws.Cells.Item(2, 4).Value = "合作金庫商業銀行臺南市北區曲門路"  # D2 (owner)
$ws.Cells.Item(2, 5).Value = 930000                 # E2 (total)
$ws.Cells.Item(2, 6).Value = "89年03月29日"          # F2 (register_date)
$ws.Cells.Item(2, 7).Value = "信用貸款"              # G2 (register_reason)
$ws.Cells.Item(2, 8).Value = "debt"                 # H2 (property_category, new string)
$ws.Cells.Item(2, 9).Value = "normal"                # I2 (category)
$ws.Cells.Item(2, 10).Value = "2012-05-01"           # J2 (date)
$ws.Cells.Item(2, 11).Value = "陳亭妃"               # K2 (legislator_name)
$ws.Cells.Item(2, 12).Value = 1708                   # L2 (legislator_id)
$ws.Cells.Item(2, 13).Value = "tmpb2a21"             # M2 (source_file)
$ws.Cells.Item(2, 14).Value = 83                     # N2 (index)

# Match row 2's plain style on the newly added cells (H2:N2).
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
